$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new sample rows right after the header row. ---
# This pushes the existing data rows (old rows 2-21) down to rows 5-24.
$ws.Range("A2:A4").EntireRow.Insert()
# Newly inserted rows pick up the header row's formatting; strip it so the
# new data rows look like the rest of the plain data rows.
$ws.Range("A2:H4").ClearFormats()

$newTopRows = @(
    @(-0.8363723754882812, 5.859383583068848, 2.452773094177246, 0.01617096064405302, -0.001510194632121166, 0.006719517832001004),
    @(-0.9548721313476562, 5.95263671875, 2.935124397277832, 0.009587190579622892, -0.01432139695518542, 0.07086037078665362),
    @(-1.523673057556152, 5.953047752380371, 3.276031017303467, 0.0104356142692267, 0.1928298026323316, 0.05640322466691333)
)

for ($i = 0; $i -lt $newTopRows.Count; $i++) {
    $r = 2 + $i
    $vals = $newTopRows[$i]
    $ws.Cells.Item($r, 2).Value = "falling"
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($r, 3 + $j).Value = $vals[$j]
    }
}

# --- Append 7 new sample rows after the (now shifted) existing data. ---
# Existing data now ends at row 24, so new rows go at 25-31.
$newBottomRows = @(
    @(1.72843074798584, 5.596723556518555, -0.66290283203125, -0.168751522898674, 0.04469497253497444, -0.09234245866537093),
    @(1.581844329833984, 5.306270599365234, -0.953785240650177, -0.1727900256713232, -0.1038810287912686, 0.06023810141616398),
    @(1.586828231811523, 5.404983997344971, -0.8601570129394531, -0.09510832776625952, -0.07622240483760839, 0.01844473597076204),
    @(1.426663398742676, 5.44196891784668, -0.7858069539070129, -0.0361937656998634, 0.0348193198442459, -0.0697913542389869),
    @(1.610628128051758, 5.431691646575928, -0.8632726669311523, -0.03700825323661167, 0.04744386838542084, 0.007177666657500789),
    @(1.379239082336426, 5.286327838897705, -0.781875729560852, -0.0347344755298561, 0.06149377011590525, 0.0545706277092296),
    @(1.496992111206055, 5.291580200195312, -0.7745996713638306, -0.01959859269360693, 0.00315613796313606, 0.01527163075904057)
)

for ($i = 0; $i -lt $newBottomRows.Count; $i++) {
    $r = 25 + $i
    $vals = $newBottomRows[$i]
    $ws.Cells.Item($r, 2).Value = "falling"
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($r, 3 + $j).Value = $vals[$j]
    }
}

# --- Recompute the timestamp column (A) for every data row. ---
# Timestamps are a simple 100 ms sequence starting at 0 on row 2.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 31) { $lastRow = 31 }
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2) * 100
}

$ws.Range("A1").Select()
